$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-13 (20 columns, A-T) reflecting updated
# NATMI LR-pair results (Ccl5-Cxcr3) with FAPs cluster added and split
# "Target cluster" values (M1/M2) per Dr Hou's advice.
$data = @(
  @("ECs", "Ccl5", "Cxcr3", "M1", 2, 1, 0.6472015, 1.294403, 0.02102898872844592, 0.01428114079438904, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 0.4638486678485, 2.783092007091, 0.00478663479322701, 0.003250684390777549),
  @("ECs", "Ccl5", "Cxcr3", "M2", 2, 1, 0.6472015, 1.294403, 0.02102898872844592, 0.01428114079438904, 3, 1, 2.431954666666666, 7.295864, 0.7723792211295388, 0.7723792211295388, 1.573964708198667, 9.443788249192, 0.01624235393521891, 0.01103045640361149),
  @("FAPs", "Ccl5", "Cxcr3", "M1", 2, 0.6666666666666666, 0.8756340000000001, 2.626902, 0.02845125902248994, 0.02898259453590742, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 0.6275660121660002, 5.648094109494001, 0.006476097738544395, 0.006597040741950017),
  @("FAPs", "Ccl5", "Cxcr3", "M2", 2, 0.6666666666666666, 0.8756340000000001, 2.626902, 0.02845125902248994, 0.02898259453590742, 3, 1, 2.431954666666666, 7.295864, 0.7723792211295388, 0.7723792211295388, 2.129502192592, 19.165519733328, 0.02197516128394554, 0.0223855537939574),
  @("M1", "Ccl5", "Cxcr3", "M1", 3, 1, 3.540639, 10.621917, 0.1150430856889177, 0.1171915486778959, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 2.537572430661001, 22.838151875949, 0.02618619676817264, 0.02667523158709822),
  @("M1", "Ccl5", "Cxcr3", "M2", 3, 1, 3.540639, 10.621917, 0.1150430856889177, 0.1171915486778959, 3, 1, 2.431954666666666, 7.295864, 0.7723792211295388, 0.7723792211295388, 8.610673539032, 77.496061851288, 0.08885688892074502, 0.09051631709079766),
  @("M2", "Ccl5", "Cxcr3", "M1", 3, 1, 23.42392733333334, 70.271782, 0.761094502822696, 0.7753081633885378, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 16.78790529587267, 151.091147662854, 0.1732409235265285, 0.1764762480151257),
  @("M2", "Ccl5", "Cxcr3", "M2", 3, 1, 23.42392733333334, 70.271782, 0.761094502822696, 0.7753081633885378, 3, 1, 2.431954666666666, 7.295864, 0.7723792211295388, 0.7723792211295388, 56.96592938996089, 512.693364509648, 0.5878535792961675, 0.5988319153734121),
  @("Neutro", "Ccl5", "Cxcr3", "M1", 3, 1, 1.243757666666667, 3.731273, 0.04041240008444284, 0.04116711337605242, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 0.8913998759423334, 8.022598883481, 0.009198701983245568, 0.00937049041050563),
  @("Neutro", "Ccl5", "Cxcr3", "M2", 3, 1, 1.243757666666667, 3.731273, 0.04041240008444284, 0.04116711337605242, 3, 1, 2.431954666666666, 7.295864, 0.7723792211295388, 0.7723792211295388, 3.024762261652444, 27.222860354872, 0.03121369810119727, 0.03179662296554678),
  @("sCs", "Ccl5", "Cxcr3", "M1", 2, 1, 1.045475, 2.09095, 0.03396976365300761, 0.02306943922721731, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 0.7492908870250002, 4.495745322150001, 0.007732224060743074, 0.005251083725003974),
  @("sCs", "Ccl5", "Cxcr3", "M2", 2, 1, 1.045475, 2.09095, 0.03396976365300761, 0.02306943922721731, 3, 1, 2.431954666666666, 7.295864, 0.7723792211295388, 0.7723792211295388, 2.542547805133334, 15.2552868308, 0.02623753959226453, 0.01781835550221334),
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}
